# quarterly.xlsx update: drop oldest quarter (1399/06), shift the remaining
# quarterly columns one step to the left, and append the newest quarter
# (1401/12) together with its reported figures - both for the headers
# (row 8 / row 24) and for every data row that reports quarterly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L","M","N")

# ---- Row 8 / Row 24 quarter-period headers ----
$labels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $labels[$i]
    $ws.Range($cols[$i] + "24").Value = $labels[$i]
}

# ---- Data rows: shift each quarterly figure one column left, append new quarter value ----
$rowData = @{}
$rowData[10] = @(177797,249893,126861,145489,198391,224465,163654,198389,159376,260612)
$rowData[13] = @(467,1118,182,383,344,1672,1826,486,1739,6097)
$rowData[14] = @(-93,610,266,275,888,351,1824,1184,1318,1645)
$rowData[15] = @(1097,-160,274,318,340,2910,1033,190,1207,1955)
$rowData[16] = @(659,726,890,891,862,862,961,1119,1631,2088)
$rowData[17] = @(21860,27242,22086,51199,14683,44188,47094,77933,48275,60461)
$rowData[19] = @(38978,92963,45696,76564,45120,86836,40750,139274,93727,43967)
$rowData[20] = @(240765,372392,196255,275119,260628,361284,257142,418575,307273,376825)
$rowData[26] = @(120,107,103,102,107,107,105,104,102,108)
$rowData[27] = @(623,637,617,637,638,640,666,634,631,730)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
